# Use a default parser if none exists
#
# The "Number" column header in the contacts fixture is renamed to "Phone".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Phone"
